$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("D8").Value = 44831
$ws.Range("J8").Value = 300

# Row 9
$ws.Range("D9").Value = 44831
$ws.Range("J9").Value = 200

# Row 10
$ws.Range("D10").Value = 44764
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 700
$ws.Range("L10").Value = 800
$ws.Range("M10").Value = 750
$ws.Range("P10").Value = 750

# Row 11
$ws.Range("D11").Value = 44764
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 150
$ws.Range("K11").Value = 600
$ws.Range("L11").Value = 600
$ws.Range("M11").Value = 600
$ws.Range("P11").Value = 600

# Row 12
$ws.Range("D12").Value = 44624
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 650
$ws.Range("L12").Value = 700
$ws.Range("M12").Value = 675
$ws.Range("P12").Value = 675

# Row 13
$ws.Range("D13").Value = 44804
$ws.Range("K13").Value = 750
$ws.Range("L13").Value = 850
$ws.Range("M13").Value = 800
$ws.Range("P13").Value = 800

# Row 14
$ws.Range("D14").Value = 44804
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 650
$ws.Range("L14").Value = 650
$ws.Range("M14").Value = 650
$ws.Range("P14").Value = 650

# Row 15
$ws.Range("D15").Value = 44761
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 700
$ws.Range("L15").Value = 800
$ws.Range("M15").Value = 750
$ws.Range("P15").Value = 750

# Row 16
$ws.Range("D16").Value = 44761
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 150
$ws.Range("L16").Value = 600
$ws.Range("M16").Value = 600
$ws.Range("P16").Value = 600

# Row 17
$ws.Range("D17").Value = 44608
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 650
$ws.Range("M17").Value = 625
$ws.Range("P17").Value = 625

$wb.Save()
